$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 (605179, one report row) gets refreshed to a newer report (2019 annual
# report superseding the previous 2020 Q3 snapshot). Update every field that
# changed per the bot's data refresh.

$ws.Range("G6").Value = "2020-12-15 00:00:00"
$ws.Range("H6").Value = "2019-12-31 00:00:00"
$ws.Range("I6").Value = 0.51
$ws.Range("J6").Value = 0.47
$ws.Range("K6").Value = 1997112513.68
$ws.Range("L6").Value = 173994761.11
$ws.Range("M6").Value = 30.29
$ws.Range("N6").Value = 13.797705172
$ws.Range("O6").Value = 10.5094616819
$ws.Range("P6").Value = 1.944785347324
$ws.Range("Q6").Value = 0.876361199941
$ws.Range("R6").Value = 39.2799002448
$ws.Range("S6").Value = 3.792
$ws.Range("T6").Value = 27.9966

# ISNEW (AB6) and DATAYEAR (AE6) hold digit-only text in the source data
# ("0", "2019"), not numbers. Force text formatting before writing so the
# engine keeps them as strings instead of auto-coercing to numeric cells,
# then drop the number format back to General so no stray style survives.
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "0"
$ws.Range("AB6").Style = "Normal"

$ws.Range("AC6").Value = "2019Q4"
$ws.Range("AD6").Value = "2019年 年报"

$ws.Range("AE6").NumberFormat = "@"
$ws.Range("AE6").Value = "2019"
$ws.Range("AE6").Style = "Normal"

$ws.Range("AF6").Value = "年报"
$ws.Range("AG6").Value = "2020-12-08 07:26:10"
